$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows 10-12: give the "No./Marking/Total" labels the mtitleStyle (like A9/row header) ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# --- Update the right/wrong/not-attempt/max counts now that the student has answered ---
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 28

# --- Marking scheme: +4 for correct, -1 for wrong (C11 becomes a real number, not text) ---
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Totals: 13*4 = 52 right-marks, 3*-1 = -3 wrong-marks, final score 49/112 ---
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "49/112"

# --- The sheet used to show three Student/Correct answer pairs (A/B, D/E, G/H); now only
#     A/B stays fully populated, D/E only for the first two questions, G/H is dropped ---
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# D16 / D18 now carry a student answer too, styled like a correct answer (correctStyle)
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"
$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# --- Fill in the student's answers (column A) for each question, styled correct (green) or
#     incorrect (red) depending on whether they match the "Correct Ans" column (B) ---
$correctCells = @{
  "A18" = "Option B"
  "A21" = "Option C"
  "A22" = "Option D"
  "A23" = "Option D"
  "A25" = "Option A"
  "A27" = "Option A"
  "A30" = "Option B"
  "A31" = "Option D"
  "A32" = "Option C"
  "A38" = "Option A"
  "A39" = "Option D"
}
foreach ($ref in $correctCells.Keys) {
  $ws.Range("B10").Copy()
  $ws.Range($ref).PasteSpecial(-4122)
  $ws.Range($ref).Value = $correctCells[$ref]
}

$incorrectCells = @{
  "A28" = "Option B"
  "A35" = "Option B"
  "A37" = "Option C"
}
foreach ($ref in $incorrectCells.Keys) {
  $ws.Range("C10").Copy()
  $ws.Range($ref).PasteSpecial(-4122)
  $ws.Range($ref).Value = $incorrectCells[$ref]
}
